$d = $word.ActiveDocument

function Insert-ItalicParagraphAfter($searchText, $newText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $pos = $r.Start + 1
    $nr = $d.Range($pos, $pos)
    $nr.Text = $newText
    $nr.Font.Italic = 1
}

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2. Add English translation of "Objetivos" text
Insert-ItalicParagraphAfter "Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais." "Provide the student with the basics of digital devices and their applications with an emphasis on microcontrollers and digital signal processors."

# 3. Add English translation of "Programa resumido" text
Insert-ItalicParagraphAfter "Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle." "Digital circuits. Microprocessors and microcontrollers. Programming of data acquisition systems and control algorithms."

# 4. Add English translation of "Programa" text
Insert-ItalicParagraphAfter "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA." "Numerical bases. Binary arithmetic. Logical functions. Boolean algebra. Minimization. Combinatorial circuits. flip-flops. Accountants and Accountants Design. Introduction to sequential circuits. Microprocessors. Microcontrollers and embedded systems. Communication interfaces. Low-level and high-level programming language in real-time computing. Development of digital command protocols. Project with programmable devices: microcontrollers and digital signal processors. Programming of FPGA devices."

# 5. Replace the two requisito lines with a single updated requirement
$reqSearch = "LOM3206 -  Eletrônica  (Requisito)`vLOM3221 -  Laboratório de Eletrônica  (Requisito)"
$reqReplace = "LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)"
$d.Content.Find.Execute($reqSearch, $false, $false, $false, $false, $false, $true, 1, $false, $reqReplace, 2)
